$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column D width: 13 -> 12
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666

# Row 2
$ws.Range("D2").Value = 52.66
$ws.Range("E2").Value = -52.66

# Row 3
$ws.Range("C3").Value = 13723.34
$ws.Range("D3").Value = 331.75
$ws.Range("E3").Value = 13391.59
$ws.Range("F3").Value = 0.02417414419521778

# Row 4
$ws.Range("C4").Value = 13723.34
$ws.Range("D4").Value = 384.41
$ws.Range("E4").Value = 13338.93
$ws.Range("F4").Value = 0.02801140247199297
